$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 54: new journal entry (Implémentation Placement multiple) ---
$ws.Range("A54").Value = 44697
$ws.Range("B54").Value = 0.5625
$ws.Range("C54").Value = 0.60416666666666663

# --- Row 55: new journal entry (Entretien avec chef de projet) ---
$ws.Range("A55").Value = 44697
$ws.Range("B55").Value = 0.60416666666666663
$ws.Range("C55").Value = 0.62847222222222221

# Write the E55 description first so the shared-string table gets the two
# new strings in the same order the authored workbook uses (95 then 96).
$ws.Range("E55").Value = "Entretien avec chef de projet"
$ws.Range("E54").Value = "Implémentation Placement multiple"

# Move the active selection the way the author's session ended up
# (was E53 before this edit).
$ws.Range("F54").Select()
